# "improve timevis filters and views"
#
# Summary of changes applied:
#  - Schedule sheet becomes the active tab (was Status); its frozen pane
#    scrolls back to the top (A2) and the selection moves to I11.
#  - Schedule sheet column widths for B:E are widened.
#  - Schedule sheet header row (row 1) height reverts to the default (no
#    longer an explicit 30pt row).
#  - A bunch of "3 to 6" (month range) values in column I become the
#    numeric estimate 4.5; a bunch of "6" values become "3" or "completed".
#  - One schedule end-date (G51) moves out a year.
#  - Row 67 gets an explicit (slightly shorter) custom height.

$wb = $excel.ActiveWorkbook
$schedule = $wb.Worksheets.Item("Schedule")
$status = $wb.Worksheets.Item("Status")

# --- Column widths (Schedule!B:E) ---------------------------------------
# The runtime's ColumnWidth COM property only resolves to 1/6-character
# granularity, so these are the closest achievable widths to the target
# OOXML widths of 36.85546875 / 11.42578125 / 27 / 28.28515625.
$schedule.Columns.Item(2).ColumnWidth = 36.0
$schedule.Columns.Item(3).ColumnWidth = 10.666666666666666
$schedule.Columns.Item(4).ColumnWidth = 26.166666666666668
$schedule.Columns.Item(5).ColumnWidth = 27.5

# --- Row 1 height reverts to default (remove explicit 30pt height) -----
$schedule.Rows.Item(1).AutoFit() | Out-Null

# --- Column I values: "3 to 6" -> 4.5 -----------------------------------
$rows45 = @(10, 26, 27, 28, 29, 30, 31, 32, 35, 49)
foreach ($r in $rows45) {
    $schedule.Cells.Item($r, 9).Value = 4.5
}

# --- G51 end date moves from 2020-02-19 (43880) to 2021-02-19 (44246) --
$schedule.Cells.Item(51, 7).Value = 44246

# --- Column I values: 6 -> 3 --------------------------------------------
$rows3 = @(51, 59, 60, 63, 64, 67, 68)
foreach ($r in $rows3) {
    $schedule.Cells.Item($r, 9).Value = 3
}

# --- Column I values: 6 -> "completed" ----------------------------------
$rowsCompleted = @(61, 62, 65, 66)
foreach ($r in $rowsCompleted) {
    $schedule.Cells.Item($r, 9).Value = "completed"
}

# --- Row 67 gets a slightly shorter, explicit custom height ------------
$schedule.Rows.Item(67).RowHeight = 29.25

# --- Sheet view / active tab changes ------------------------------------
# Activating Schedule sets its tabSelected + the workbook's activeTab
# (and clears Status's tabSelected since it was previously active), and
# resets the frozen pane's scroll position back to A2.
$schedule.Activate() | Out-Null
$schedule.Range("I11").Select() | Out-Null
